$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price data as literal text (e.g. "1.000", "0.9999") in the
# source workbook. Force Text format on the whole data range first so that
# writing number-looking strings does not get reinterpreted as numeric values.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '31.107.61'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '1.956.36'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '246.81'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '0.4874'
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '0.2967'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('D10').Value = '0.06850'
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('D11').Value = '19.09'
$ws.Range('E11').Value = '  -2.86%  '
$ws.Range('D12').Value = '106.43'
$ws.Range('E12').Value = '  -5.36%  '
$ws.Range('D13').Value = '1.937.40'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '0.07743'
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('D15').Value = '5.416'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').Value = '0.7138'
$ws.Range('E16').Value = '  +4.26%  '
$ws.Range('D17').Value = '284.57'
$ws.Range('E17').Value = '  -5.17%  '
$ws.Range('D18').Value = '31.119.11'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').Value = '0.000007773'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = '13.23'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').Value = '0.4935'
$ws.Range('E21').Value = '  +8.88%  '
$ws.Range('D22').Value = '2.194.29'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').Value = '0.9994'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '5.549'
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').Value = '0.9997'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').Value = '6.589'
$ws.Range('E26').Value = '  +0.96%  '
$ws.Range('D27').Value = '9.933'
$ws.Range('E27').Value = '  +4.10%  '
$ws.Range('D28').Value = '169.28'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('D29').Value = '19.93'
$ws.Range('E29').Value = '  -3.34%  '
$ws.Range('D30').Value = '2.195'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('D31').Value = '0.1050'
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').Value = '1.439'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').Value = '4.758'
$ws.Range('E33').Value = '  +15.87%  '
$ws.Range('D34').Value = '4.500'
$ws.Range('E34').Value = '  +7.68%  '
$ws.Range('D35').Value = '0.05010'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').Value = '0.7657'
$ws.Range('E36').Value = '  +2.74%  '
$ws.Range('D37').Value = '1.168'
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').Value = '2.731'
$ws.Range('E39').Value = '  +0.81%  '
$ws.Range('D40').Value = '2.715'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '2.156'
$ws.Range('E41').Value = '  +5.31%  '
$ws.Range('D42').Value = '6.451'
$ws.Range('E42').Value = '  +9.57%  '
$ws.Range('D43').Value = '109.76'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('D44').Value = '0.4472'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '0.8819'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = '72.77'
$ws.Range('E46').Value = '  +3.84%  '
$ws.Range('D47').Value = '0.9993'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').Value = '7.518'
$ws.Range('E48').Value = '  +2.77%  '
$ws.Range('D49').Value = '994.92'
$ws.Range('E49').Value = '  +17.24%  '
$ws.Range('D50').Value = '0.1271'
$ws.Range('E50').Value = '  +2.49%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '0.2631'
$ws.Range('E51').Value = '  +3.35%  '
